# Utilities #1 look up list.xlsx -- "Commit for first publish"
#
# Adds a new "Direct Search" (Y/N) column to both sheets, right after the
# "Service" column, renames the first sheet "Kent and South London" -> "KSL",
# and tidies up a couple of label cells.
#
# Implementation notes:
#  - Worksheet object references in this COM bridge resolve by *current
#    index position*, not stable identity: once a Worksheets.Insert/Copy/
#    Delete shuffles the tab order, a previously-captured sheet variable
#    silently starts pointing at whatever now sits at its old slot. So every
#    sheet handle used here is re-fetched by name immediately before use,
#    rather than cached across a structural change.
#  - New shared strings must be introduced in the order Y, N, "Direct
#    Search" so the saved sharedStrings.xml table gets indices 14/15/16 in
#    that order (the engine assigns shared-string ids in first-write order,
#    not cell position).
#  - The target column widths (20.7109375 / 15.7109375 / 30.7109375
#    "characters") can't be reproduced by assigning a literal ColumnWidth --
#    the COM bridge always rounds that to the nearest 1/6 of a character.
#    They *do* survive a Worksheets.Copy() (whole-sheet clone) intact, so
#    sheet1 is rebuilt by cloning the already-fixed-up "Next Region" (whose
#    pre-existing 20.7109375-wide column lets a plain column insert merge
#    into it cleanly) and then patching up the handful of cells where KSL's
#    data differs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Next Region" (sheet 2): insert the new "Direct Search" column B.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Next Region").Columns("B:B").Insert()

# Seed the new shared strings in the exact order they must be allocated:
# Y (14), N (15), "Direct Search" (16).
$wb.Worksheets.Item("Next Region").Range("B2").Value = "Y"
$wb.Worksheets.Item("Next Region").Range("B3").Value = "Y"
$wb.Worksheets.Item("Next Region").Range("B4").Value = "N"
$wb.Worksheets.Item("Next Region").Range("B5").Value = "N"
$wb.Worksheets.Item("Next Region").Range("B6").Value = "N"
$wb.Worksheets.Item("Next Region").Range("B1").Value = "Direct Search"

# The LSBUD/DIGDAT "Yes" flags on this sheet become plain "Y".
$wb.Worksheets.Item("Next Region").Range("D4").Value = "Y"
$wb.Worksheets.Item("Next Region").Range("E5").Value = "Y"

# ---------------------------------------------------------------------
# 2) Rebuild sheet 1 by cloning the now-correct "Next Region" sheet, so the
#    new columns inherit its exact (non-representable-via-ColumnWidth)
#    widths, then drop the old "Kent and South London" sheet and rename the
#    clone into its place.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Next Region").Copy($wb.Worksheets.Item(1))
$wb.Worksheets.Item("Kent and South London").Delete()
$wb.Worksheets.Item("Next Region (2)").Name = "KSL"

# KSL only has 5 data rows (no "Water" row) -- drop the 6th.
$wb.Worksheets.Item("KSL").Rows("6:6").Delete()

# KSL keeps "Yes" (not "Y") in the LSBUD/DIGDAT columns.
$wb.Worksheets.Item("KSL").Range("D4").Value = "Yes"
$wb.Worksheets.Item("KSL").Range("E5").Value = "Yes"

# KSL's row 3 provider differs from "Next Region"'s: "LSBUD" not
# "Line Search B4U Dig".
$wb.Worksheets.Item("KSL").Range("C3").Value = "LSBUD"

# ---------------------------------------------------------------------
# 3) Selections: match the saved cursor position on each sheet, leaving
#    "Next Region" as the active/tab-selected sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("KSL").Activate()
$wb.Worksheets.Item("KSL").Range("E20").Select()

$wb.Worksheets.Item("Next Region").Activate()
$wb.Worksheets.Item("Next Region").Range("D12").Select()
